# Update "想去人数" (interested-attendee count, column F) figures to the latest
# scraped values, matching the gh-pages data refresh at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 226   # was 179
$ws.Range("F3").Value = 1058   # was 1053
$ws.Range("F4").Value = 9383   # was 9316
$ws.Range("F5").Value = 197   # was 196
$ws.Range("F6").Value = 69   # was 68
$ws.Range("F7").Value = 6448   # was 6421
$ws.Range("F8").Value = 625   # was 622
$ws.Range("F9").Value = 71   # was 67
$ws.Range("F10").Value = 9857   # was 9812
$ws.Range("F11").Value = 11241   # was 11167
$ws.Range("F12").Value = 1242   # was 1236
$ws.Range("F13").Value = 1159   # was 1152
$ws.Range("F14").Value = 4946   # was 4929
$ws.Range("F15").Value = 798   # was 795
$ws.Range("F16").Value = 464   # was 456
$ws.Range("F17").Value = 96   # was 95
$ws.Range("F18").Value = 333   # was 332
$ws.Range("F19").Value = 178   # was 177
$ws.Range("F20").Value = 1341   # was 1336
$ws.Range("F21").Value = 252   # was 246
$ws.Range("F22").Value = 1858   # was 1851
$ws.Range("F23").Value = 887   # was 886
$ws.Range("F24").Value = 1255   # was 1247
$ws.Range("F25").Value = 857   # was 856
$ws.Range("F27").Value = 2044   # was 2033
$ws.Range("F28").Value = 431   # was 425
$ws.Range("F29").Value = 626   # was 622
$ws.Range("F30").Value = 2677   # was 2665
$ws.Range("F31").Value = 186   # was 185
$ws.Range("F32").Value = 1769   # was 1750
$ws.Range("F35").Value = 61   # was 56
$ws.Range("F36").Value = 919   # was 914
$ws.Range("F37").Value = 588   # was 589
$ws.Range("F38").Value = 27   # was 24
$ws.Range("F39").Value = 3341   # was 3321
$ws.Range("F40").Value = 236   # was 234
$ws.Range("F41").Value = 84   # was 82
$ws.Range("F42").Value = 515   # was 512
$ws.Range("F43").Value = 579   # was 578
$ws.Range("F46").Value = 240   # was 238
$ws.Range("F48").Value = 4214   # was 4210
$ws.Range("F49").Value = 43   # was 32

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 11   # was 10
$ws.Range("F9").Value = 4   # was 0
$ws.Range("F26").Value = 44   # was 43

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5948   # was 5933

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 226   # was 179
$ws.Range("F3").Value = 1058   # was 1053
$ws.Range("F4").Value = 9383   # was 9316
$ws.Range("F5").Value = 197   # was 196
$ws.Range("F6").Value = 69   # was 68
$ws.Range("F7").Value = 11   # was 10
$ws.Range("F8").Value = 6448   # was 6421
$ws.Range("F9").Value = 625   # was 622
$ws.Range("F10").Value = 9857   # was 9812
$ws.Range("F11").Value = 11241   # was 11168
$ws.Range("F13").Value = 1159   # was 1152
$ws.Range("F14").Value = 4946   # was 4929
$ws.Range("F15").Value = 798   # was 795
$ws.Range("F16").Value = 464   # was 456
$ws.Range("F17").Value = 96   # was 95
$ws.Range("F18").Value = 333   # was 332
$ws.Range("F20").Value = 178   # was 177
$ws.Range("F21").Value = 1341   # was 1336
$ws.Range("F22").Value = 252   # was 246
$ws.Range("F23").Value = 1858   # was 1851
$ws.Range("F24").Value = 857   # was 856
$ws.Range("F26").Value = 2044   # was 2033
$ws.Range("F27").Value = 431   # was 425
$ws.Range("F28").Value = 626   # was 622
$ws.Range("F29").Value = 2677   # was 2665
$ws.Range("F30").Value = 186   # was 185
$ws.Range("F31").Value = 1769   # was 1750
$ws.Range("F39").Value = 61   # was 56
$ws.Range("F40").Value = 919   # was 914
$ws.Range("F41").Value = 588   # was 589
$ws.Range("F42").Value = 27   # was 24
$ws.Range("F43").Value = 44   # was 43
$ws.Range("F44").Value = 237   # was 234
$ws.Range("F45").Value = 579   # was 578
$ws.Range("F47").Value = 240   # was 238
$ws.Range("F49").Value = 4214   # was 4210
